# Insert a new data row for "Macroferia Regional de Talca - Repollo" right before
# the current row 392 (the other rows below shift down by one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(392).Insert()

$ws.Cells.Item(392, 1).Value = 5
$ws.Cells.Item(392, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(392, 3).Value = "Maule"
$ws.Cells.Item(392, 4).Value = 44931
$ws.Cells.Item(392, 5).Value = 7
$ws.Cells.Item(392, 6).Value = 100112006
$ws.Cells.Item(392, 7).Value = "Repollo"
$ws.Cells.Item(392, 8).Value = "Crespo record"
$ws.Cells.Item(392, 9).Value = "Primera"
$ws.Cells.Item(392, 10).Value = 3000
$ws.Cells.Item(392, 11).Value = 1000
$ws.Cells.Item(392, 12).Value = 1000
$ws.Cells.Item(392, 13).Value = 1000
$ws.Cells.Item(392, 14).Value = "`$/unidad"
$ws.Cells.Item(392, 15).Value = "Región del Maule"
$ws.Cells.Item(392, 16).Value = 1000
$ws.Cells.Item(392, 17).Value = 1
$ws.Cells.Item(392, 18).Value = "Hortaliza"
